$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 54
$ws1.Range("F5").Value = 206
$ws1.Range("F6").Value = 1131
$ws1.Range("F7").Value = 1079
$ws1.Range("F8").Value = 8380
$ws1.Range("F11").Value = 6964
$ws1.Range("F13").Value = 305
$ws1.Range("F14").Value = 5118
$ws1.Range("F15").Value = 5118
$ws1.Range("F17").Value = 62
$ws1.Range("F18").Value = 5621
$ws1.Range("F19").Value = 5621
$ws1.Range("F20").Value = 1082
$ws1.Range("F21").Value = 347
$ws1.Range("F24").Value = 486
$ws1.Range("F26").Value = 259
$ws1.Range("F28").Value = 9417
$ws1.Range("F30").Value = 1731
$ws1.Range("F31").Value = 1110
$ws1.Range("F32").Value = 41
$ws1.Range("F33").Value = 1912
$ws1.Range("F38").Value = 1919
$ws1.Range("F40").Value = 1241
$ws1.Range("F41").Value = 55
$ws1.Range("F42").Value = 4911
$ws1.Range("F44").Value = 1168
$ws1.Range("F45").Value = 533
$ws1.Range("F48").Value = 1050
$ws1.Range("F50").Value = 1287

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 4
$ws2.Range("F10").Value = 35
$ws2.Range("F11").Value = 184

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 54
$ws4.Range("F6").Value = 206
$ws4.Range("F8").Value = 1131
$ws4.Range("F9").Value = 1079
$ws4.Range("F10").Value = 8380
$ws4.Range("F13").Value = 6964
$ws4.Range("F15").Value = 305
$ws4.Range("F18").Value = 5118
$ws4.Range("F19").Value = 5118
$ws4.Range("F21").Value = 5621
$ws4.Range("F22").Value = 5621
$ws4.Range("F23").Value = 1082
$ws4.Range("F24").Value = 347
$ws4.Range("F26").Value = 486
$ws4.Range("F28").Value = 259
$ws4.Range("F30").Value = 184
$ws4.Range("F31").Value = 9417
$ws4.Range("F33").Value = 1731
$ws4.Range("F34").Value = 1110
$ws4.Range("F35").Value = 41
$ws4.Range("F36").Value = 1912
$ws4.Range("F40").Value = 1919
$ws4.Range("F42").Value = 1241
$ws4.Range("F43").Value = 4911
$ws4.Range("F45").Value = 1168
$ws4.Range("F46").Value = 533
$ws4.Range("F49").Value = 1050
$ws4.Range("F51").Value = 1287
